$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.341.54'
$ws.Range("E2").Value = '  +0.64%  '

$ws.Range("D3").Value = '1.880.53'
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("E4").Value = '  -0.64%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.96'
$ws.Range("E5").Value = '  -2.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.683'
$ws.Range("E6").Value = '  -1.87%  '

$ws.Range("E7").Value = '  -0.63%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.51'
$ws.Range("E8").Value = '  +4.77%  '

$ws.Range("E9").Value = '  +1.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.37'
$ws.Range("E10").Value = '  +2.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0745'
$ws.Range("E11").Value = '  -1.20%  '

$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.56'
$ws.Range("E13").Value = '  +4.10%  '

$ws.Range("D14").Value = '2.153.54'
$ws.Range("E14").Value = '  -0.94%  '

$ws.Range("E15").Value = '  +5.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.94'
$ws.Range("E16").Value = '  -0.41%  '

$ws.Range("D17").Value = '1.893.11'
$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").Value = '35.370.40'
$ws.Range("E18").Value = '  +0.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.61'
$ws.Range("E19").Value = '  -1.03%  '

$ws.Range("E20").Value = '  -1.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '244.72'
$ws.Range("E21").Value = '  -3.00%  '

$ws.Range("E22").Value = '  -1.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.09'
$ws.Range("E23").Value = '  +1.22%  '

$ws.Range("E24").Value = '  +10.72%  '

$ws.Range("E25").Value = '  -0.57%  '

$ws.Range("E26").Value = '  -4.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.63'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.74'
$ws.Range("E28").Value = '  +1.95%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.31'
$ws.Range("E29").Value = '  -0.53%  '

$ws.Range("E30").Value = '  -0.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.30'
$ws.Range("E31").Value = '  -0.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0595'
$ws.Range("E32").Value = '  +0.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.19'
$ws.Range("E33").Value = '  -1.28%  '

$ws.Range("E34").Value = '  -0.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.80'
$ws.Range("E35").Value = '  -9.91%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.43'
$ws.Range("E36").Value = '  -11.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.854'
$ws.Range("E37").Value = '  +0.84%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.96'
$ws.Range("E38").Value = '  -3.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0728'
$ws.Range("E39").Value = '  +9.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.38'
$ws.Range("E40").Value = '  -0.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0218'
$ws.Range("E41").Value = '  +2.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.15'
$ws.Range("E42").Value = '  -1.73%  '

$ws.Range("E43").Value = '  -2.26%  '

$ws.Range("E44").Value = '  +0.49%  '

$ws.Range("D45").Value = '1.309.09'
$ws.Range("E45").Value = '  +0.52%  '

$ws.Range("E46").Value = '  +4.25%  '

$ws.Range("E47").Value = '  -1.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.74'
$ws.Range("E48").Value = '  -0.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.91'
$ws.Range("E49").Value = '  -1.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.31'
$ws.Range("E50").Value = '  -4.15%  '

$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '42.14'
$ws.Range("E51").Value = '  -1.48%  '
